$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting the existing rows 19-37 down to 20-38
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 45097
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100108
$ws.Range("H19").Value = "Tropicales y subtropicales"
$ws.Range("I19").Value = 100108001
$ws.Range("J19").Value = "Guayaba"
$ws.Range("K19").Value = "Sin especificar"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 90
$ws.Range("N19").Value = 2600
$ws.Range("O19").Value = 2600
$ws.Range("P19").Value = 2600
$ws.Range("Q19").Value = "$/kilo"
$ws.Range("R19").Value = "Región de Arica y Parinacota"
$ws.Range("S19").Value = 2600
$ws.Range("T19").Value = 1
